$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.138.78"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.578.46"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.31"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "657.85"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.401"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "3.574.18"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.41"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.47"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "4.240.58"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "95.073.15"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "3.540.51"
$ws.Range("E18").Value = "  -3.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("E20").Value = "  -6.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.88"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.47"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.485"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "508.82"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.01"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.03"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.82"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "3.766.52"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.57"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.177"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.12"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("E37").Value = "  +17.03%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.69"
$ws.Range("E38").Value = "  +10.79%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.561"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "602.49"
$ws.Range("E40").Value = "  +6.19%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.913"
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("E44").Value = "  +7.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "35.20"
$ws.Range("E45").Value = "  +15.87%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.40"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.51"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  +0.10%  "
